# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handoff packages have moved from "In Translation" to
# "Ready for handoff", and refreshes the associated timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 16:40:43"

# Status text got longer ("In Translation" -> "Ready for handoff"), so the
# status columns widen accordingly.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 16:40:38"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-24 16:40:43"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
